$data = @(
    @("ECs", "Il6", "Il6ra", "ECs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 1.940177333333333, 5.820532, 0.09618640556857148, 0.09618640556857148, 23.08545610407467, 207.769104936672, 0.01750343562672334, 0.01750343562672334),
    @("ECs", "Il6", "Il6ra", "FAPs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 1.922884666666667, 5.768654000000001, 0.09532910277424163, 0.09532910277424164, 22.87969702710934, 205.9172732439841, 0.017347428713018, 0.017347428713018),
    @("ECs", "Il6", "Il6ra", "M2", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 15.72357133333333, 47.170714, 0.7795131832903062, 0.7795131832903062, 187.0889890210827, 1683.800901189744, 0.1418512184050491, 0.1418512184050491),
    @("ECs", "Il6", "Il6ra", "sCs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 0.5843806666666667, 1.753142, 0.02897130836688065, 0.02897130836688065, 6.953330500581334, 62.57997450523201, 0.005272028079478819, 0.005272028079478819),
    @("FAPs", "Il6", "Il6ra", "ECs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 1.940177333333333, 5.820532, 0.09618640556857148, 0.09618640556857148, 97.09370415419555, 873.84333738776, 0.07361662653583521, 0.07361662653583521),
    @("FAPs", "Il6", "Il6ra", "FAPs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 1.922884666666667, 5.768654000000001, 0.09532910277424163, 0.09532910277424164, 96.22831467019111, 866.0548320317201, 0.07296048662432437, 0.07296048662432437),
    @("FAPs", "Il6", "Il6ra", "M2", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 15.72357133333333, 47.170714, 0.7795131832903062, 0.7795131832903062, 786.8661060291688, 7081.794954262519, 0.5966033407198335, 0.5966033407198335),
    @("FAPs", "Il6", "Il6ra", "sCs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 0.5843806666666667, 1.753142, 0.02897130836688065, 0.02897130836688065, 29.24458635195111, 263.20127716756, 0.02217329960187269, 0.02217329960187269),
    @("M2", "Il6", "Il6ra", "ECs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 1.940177333333333, 5.820532, 0.09618640556857148, 0.09618640556857148, 1.567497723534222, 14.107479511808, 0.001188479680679732, 0.001188479680679732),
    @("M2", "Il6", "Il6ra", "FAPs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 1.922884666666667, 5.768654000000001, 0.09532910277424163, 0.09532910277424164, 1.553526724508445, 13.981740520576, 0.001177886843311206, 0.001177886843311206),
    @("M2", "Il6", "Il6ra", "M2", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 15.72357133333333, 47.170714, 0.7795131832903062, 0.7795131832903062, 12.70330389257956, 114.329735033216, 0.009631668567779538, 0.009631668567779538),
    @("M2", "Il6", "Il6ra", "sCs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 0.5843806666666667, 1.753142, 0.02897130836688065, 0.02897130836688065, 0.4721297115164445, 4.249167403648, 0.0003579696227675111, 0.0003579696227675111),
    @("sCs", "Il6", "Il6ra", "ECs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 1.940177333333333, 5.820532, 0.09618640556857148, 0.09618640556857148, 5.114553206462222, 46.03097885816, 0.003877863725333196, 0.003877863725333196),
    @("sCs", "Il6", "Il6ra", "FAPs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 1.922884666666667, 5.768654000000001, 0.09532910277424163, 0.09532910277424164, 5.068967546724445, 45.62070792052, 0.00384330059358805, 0.003843300593588051),
    @("sCs", "Il6", "Il6ra", "M2", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 15.72357133333333, 47.170714, 0.7795131832903062, 0.7795131832903062, 41.44932568703555, 373.04393118332, 0.03142695559764412, 0.03142695559764412),
    @("sCs", "Il6", "Il6ra", "sCs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 0.5843806666666667, 1.753142, 0.02897130836688065, 0.02897130836688065, 1.540501458884445, 13.86451312996, 0.001168011062761633, 0.001168011062761633)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowvals = $data[$r]
    for ($c = 0; $c -lt $rowvals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowvals[$c]
    }
}